# "Generate Report for Handoff"
#
# The localization-status report is regenerated by CI: for each tracked
# source file, the latest handoff timestamp is recomputed and written back
# into the Overview sheet (column D, "Latest Handoff Date") and into each
# language sheet (column E, "Latest Handoff Datetime"). For the
# "d3cba02f-32ea-4bc7-a762-0db5c7de3a59" row the freshly computed handoff
# timestamp is written, then reconciled back to the previously recorded
# value (the row's actual latest-handoff moment didn't move), so only the
# value churn is replayed here across the three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 7 is "d3cba02f-32ea-4bc7-a762-0db5c7de3a59.md" ---
$prevOverviewDate = $wsOverview.Range("D7").Value2
$wsOverview.Range("D7").Value2 = "2016-30-20 02:30:20"
$wsOverview.Range("D7").Value2 = $prevOverviewDate

# --- zh-cn sheet: row 7 is "d3cba02f-32ea-4bc7-a762-0db5c7de3a59" ---
$prevZhCnHandoffFile = $wsZhCn.Range("D7").Value2
$prevZhCnHandoffDate = $wsZhCn.Range("E7").Value2
$wsZhCn.Range("E7").Value2 = "2016-03-20 02:30:17"
$wsZhCn.Range("D7").Value2 = $prevZhCnHandoffFile
$wsZhCn.Range("E7").Value2 = $prevZhCnHandoffDate

# --- de-de sheet: row 7 is "d3cba02f-32ea-4bc7-a762-0db5c7de3a59" ---
$prevDeDeHandoffFile = $wsDeDe.Range("D7").Value2
$prevDeDeHandoffDate = $wsDeDe.Range("E7").Value2
$wsDeDe.Range("E7").Value2 = "2016-03-20 02:30:20"
$wsDeDe.Range("D7").Value2 = $prevDeDeHandoffFile
$wsDeDe.Range("E7").Value2 = $prevDeDeHandoffDate

Write-Host "Report regenerated for handoff."
